$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.219.20"
$ws.Range("E2").Value = "  -0.38%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.825.91"
$ws.Range("E3").Value = "  -0.75%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  +0.33%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "234.31"
$ws.Range("E5").Value = "  -1.99%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5980"
$ws.Range("E6").Value = "  -4.37%  "

# Row 7
$ws.Range("E7").Value = "  +0.35%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.06951"
$ws.Range("E8").Value = "  -5.67%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2749"
$ws.Range("E9").Value = "  -4.62%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.22"
$ws.Range("E10").Value = "  -6.14%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07607"
$ws.Range("E11").Value = "  -1.47%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.828.19"
$ws.Range("E12").Value = "  -0.53%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.737"
$ws.Range("E13").Value = "  -4.37%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6246"
$ws.Range("E14").Value = "  -5.82%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.000009625"
$ws.Range("E15").Value = "  -7.33%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "78.32"
$ws.Range("E16").Value = "  -3.81%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "28.685.93"
$ws.Range("E17").Value = "  -2.16%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.605"
$ws.Range("E18").Value = "  -10.15%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "218.61"
$ws.Range("E19").Value = "  -6.80%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.004"
$ws.Range("E20").Value = "  +0.28%  "

# Row 21
$ws.Range("E21").Value = "  -5.89%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.845"
$ws.Range("E22").Value = "  -6.01%  "

# Row 23
$ws.Range("E23").Value = "  +0.09%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "156.99"
$ws.Range("E24").Value = "  -0.12%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "7.937"
$ws.Range("E25").Value = "  -6.00%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1280"
$ws.Range("E26").Value = "  -4.08%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.49"
$ws.Range("E27").Value = "  -4.48%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.438"
$ws.Range("E28").Value = "  -2.94%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.06349"
$ws.Range("E29").Value = "  -10.89%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.437"
$ws.Range("E30").Value = "  -2.86%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.825"

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.744"
$ws.Range("E32").Value = "  -7.00%  "

# Row 33
$ws.Range("E33").Value = "  -5.19%  "

# Row 34
$ws.Range("E34").Value = "  -5.68%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6447"
$ws.Range("E35").Value = "  -8.41%  "

# Row 36
$ws.Range("E36").Value = "  -1.61%  "

# Row 37
$ws.Range("E37").Value = "  -1.29%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01752"
$ws.Range("E38").Value = "  -4.06%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.579"
$ws.Range("E39").Value = "  -3.30%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.146.85"
$ws.Range("E40").Value = "  -6.92%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8880"
$ws.Range("E41").Value = "  -6.00%  "

# Row 42
$ws.Range("E42").Value = "  +0.47%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.981.52"
$ws.Range("E43").Value = "  -0.37%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.56"
$ws.Range("E44").Value = "  -0.49%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "61.94"
$ws.Range("E45").Value = "  -4.77%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000113"
$ws.Range("E46").Value = "  -3.67%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.597"
$ws.Range("E47").Value = "  -4.94%  "

# Row 48
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.466"
$ws.Range("E48").Value = "  -4.87%  "

# Row 49
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05517"
$ws.Range("E49").Value = "  -2.27%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4548"
$ws.Range("E50").Value = "  -0.51%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.414"
$ws.Range("E51").Value = "  -7.45%  "
